$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item(4)

# --- Row 4: the old "datetime" combined row is cleared out (values removed,
#     formatting/styles are left untouched) ---
$ws4.Range("A4").Value = ""
$ws4.Range("B4").Value = ""
$ws4.Range("E4").Value = ""
$ws4.Range("F4").Value = ""

# --- Row 5: "_date" row now gets its column name / dtype / format filled in ---
$ws4.Range("B5").Value = "date"
$ws4.Range("E5").Value = "object"
$ws4.Range("F5").Value = "%Y-%m-%d"

# --- Row 6: "_time" row now gets its column name / dtype / format filled in ---
$ws4.Range("B6").Value = "time"
$ws4.Range("E6").Value = "object"
$ws4.Range("F6").Value = "%H:%M:%S"

# --- Row 8: the "temp" variable's template column name is now "temp" (was "value") ---
$ws4.Range("B8").Value = "temp"

# --- Row 9 (new): "flags" variable row, copy the formatting used by column A
#     of the other data rows, then fill in the values ---
$ws4.Range("A8").Copy()
$ws4.Range("A9").PasteSpecial(-4122)
$ws4.Range("A9").Value = "flags"
$ws4.Range("B9").Value = "flags"
$ws4.Range("E9").Value = "object"

# --- Rows 10-12 (new): blank rows, formatted like column A of the rest of
#     the sheet but with no content ---
$ws4.Range("A8").Copy()
$ws4.Range("A10").PasteSpecial(-4122)
$ws4.Range("A8").Copy()
$ws4.Range("A11").PasteSpecial(-4122)
$ws4.Range("A8").Copy()
$ws4.Range("A12").PasteSpecial(-4122)

# --- Update the active selection to reflect where the editor ended up ---
$ws4.Activate() | Out-Null
$ws4.Range("A6").Select() | Out-Null
